$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the "Absent" column (H) so it reflects whether a "Real"
# attendance (column E) was recorded for that day: Absent = 1 - Real.
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
